$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '72.842.20'
$ws.Range('E2').Value = '  +3.02%  '

$ws.Range('D3').Value = '3.982.14'
$ws.Range('E3').Value = '  +0.94%  '

$ws.Range('D4').Value = '''0.999'
$ws.Range('E4').Value = '  -0.20%  '

$ws.Range('D5').Value = '''591.74'
$ws.Range('E5').Value = '  +10.08%  '

$ws.Range('D6').Value = '''159.56'
$ws.Range('E6').Value = '  +7.92%  '

$ws.Range('D7').Value = '''0.686'
$ws.Range('E7').Value = '  -0.14%  '

$ws.Range('D8').Value = '''0.998'
$ws.Range('E8').Value = '  -0.21%  '

$ws.Range('D9').Value = '''0.751'
$ws.Range('E9').Value = '  +2.08%  '

$ws.Range('E10').Value = '  +2.19%  '

$ws.Range('D11').Value = '''53.94'
$ws.Range('E11').Value = '  -1.13%  '

$ws.Range('D12').Value = '''0.0000319'
$ws.Range('E12').Value = '  +0.99%  '

$ws.Range('D13').Value = '''10.93'
$ws.Range('E13').Value = '  +3.67%  '

$ws.Range('D14').Value = '4.609.29'
$ws.Range('E14').Value = '  +0.60%  '

$ws.Range('D15').Value = '3.980.56'
$ws.Range('E15').Value = '  +0.68%  '

$ws.Range('D16').Value = '''1.28'
$ws.Range('E16').Value = '  +9.93%  '

$ws.Range('D17').Value = '''14.08'
$ws.Range('E17').Value = '  +2.68%  '

$ws.Range('D18').Value = '''20.39'
$ws.Range('E18').Value = '  +0.38%  '

$ws.Range('E19').Value = '  +0.29%  '

$ws.Range('D20').Value = '72.590.09'
$ws.Range('E20').Value = '  +2.50%  '

$ws.Range('D21').Value = '''435.68'
$ws.Range('E21').Value = '  +3.10%  '

$ws.Range('D22').Value = '''4.77'
$ws.Range('E22').Value = '  +13.92%  '

$ws.Range('D23').Value = '''96.13'
$ws.Range('E23').Value = '  -0.21%  '

$ws.Range('D24').Value = '''3.43'
$ws.Range('E24').Value = '  -4.23%  '

$ws.Range('B25').Value = 'Toncoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D25').Value = '''4.48'
$ws.Range('E25').Value = '  +22.17%  '

$ws.Range('B26').Value = 'InternetComputer(DFINITY)'
$ws.Range('C26').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D26').Value = '''14.32'
$ws.Range('E26').Value = '  +1.46%  '

$ws.Range('D27').Value = '''11.19'
$ws.Range('E27').Value = '  -2.73%  '

$ws.Range('D28').Value = '''10.59'
$ws.Range('E28').Value = '  -1.28%  '

$ws.Range('E29').Value = '  +1.31%  '

$ws.Range('D30').Value = '''36.37'
$ws.Range('E30').Value = '  +0.14%  '

$ws.Range('E31').Value = '  +0.96%  '

$ws.Range('D32').Value = '''13.72'
$ws.Range('E32').Value = '  +3.00%  '

$ws.Range('E33').Value = '  +1.71%  '

$ws.Range('D34').Value = '''678.92'
$ws.Range('E34').Value = '  -1.95%  '

$ws.Range('D35').Value = '''48.27'
$ws.Range('E35').Value = '  -3.62%  '

$ws.Range('D36').Value = '''69.52'
$ws.Range('E36').Value = '  +8.30%  '

$ws.Range('D37').Value = '0.0₃0879'
$ws.Range('E37').Value = '  +8.21%  '

$ws.Range('D38').Value = '''0.435'
$ws.Range('E38').Value = '  -0.20%  '

$ws.Range('D39').Value = '''3.40'
$ws.Range('E39').Value = '  -2.25%  '

$ws.Range('E40').Value = '  -1.73%  '

$ws.Range('D41').Value = '''1.00'
$ws.Range('E41').Value = '  -0.04%  '

$ws.Range('D42').Value = '''3.35'
$ws.Range('E42').Value = '  +3.81%  '

$ws.Range('E43').Value = '  +0.08%  '

$ws.Range('D44').Value = '''10.85'
$ws.Range('E44').Value = '  +11.95%  '

$ws.Range('D45').Value = '''0.0487'
$ws.Range('E45').Value = '  +1.55%  '

$ws.Range('E46').Value = '  +1.35%  '

$ws.Range('E47').Value = '  -3.86%  '

$ws.Range('E48').Value = '  +0.43%  '

$ws.Range('E49').Value = '  +1.35%  '

$ws.Range('E50').Value = '  +4.79%  '

$ws.Range('D51').Value = '2.792.79'
$ws.Range('E51').Value = '  +11.18%  '
